# Refresh the coin price/volume snapshot to match the Sat Feb 11 2023 18:26:44 UTC
# GitHub Actions run (cryptos.xlsx symbol-list update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each of these "Price" (D) / "Volume(1h)" (E) cells holds plain text in the
# workbook (e.g. "307.97", "0.31%"), not a number. Setting NumberFormat to "@"
# (Text) before assigning the value keeps Excel from reinterpreting these
# numeric-looking / percent-looking strings as Number or Percentage values.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.31%"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.72"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.43%"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.119"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.39%"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07623"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.08%"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.619"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.10%"
# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.51%"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9088"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.40%"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1277"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "30.93%"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1813"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.91%"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09145"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.39%"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04311"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.25%"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1045"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.73%"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001254"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.16%"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005864"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.46%"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.349"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.62%"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.282"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.12%"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3314"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.49%"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.917"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.52%"
# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.35%"
# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-4.72%"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04048"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.67%"
# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "4.74%"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004062"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.01%"
# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.20%"
# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "24.87%"
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02433"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "0.40%"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05241"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.79%"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007837"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1301"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.22%"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006805"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.64%"
# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.78%"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007380"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-11.82%"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3345"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.37%"
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "8.10%"
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.25%"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1058"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1,856.51%"
# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.25%"
# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.25%"
